# Add developer name
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new unit-test rows (Method Inputs / Condition / Expected Result)
$ws.Range("C3").Value = "Ridham Sood"
$ws.Range("E7").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount`nfrom datetime import date"
$ws.Range("E8").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount"
$ws.Range("E9").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount"
$ws.Range("E10").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount"
$ws.Range("E11").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount"
$ws.Range("E12").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount"
$ws.Range("E13").Value = "import unittest`nfrom bank_account.investing_account import InvestingAccount"
$ws.Range("F7").Value = "investment = InvestingAccount(12345, 123, 1000, date(2025, 5, 1), 3.00)"
$ws.Range("F8").Value = "investment = InvestingAccount(12345, 123, 1000, date(2025, 5, 1), `"five`")`nexpected = 2.55"
$ws.Range("F9").Value = "investment = InvestingAccount(12345, 123, 1000, date(2013, 5, 1), 3.00)`nactual = investment.get_service_charge()`nexpected = 0.50"
$ws.Range("F11").Value = "investment = InvestingAccount(12345, 123, 1000, date(2018, 5, 1), 3.00)`nactual = investment.get_service_charge()`nexpected = 0.50"
$ws.Range("F10").Value = "investment = InvestingAccount(12345, 123, 1000, date(2015, 10, 2), 3.00)`nactual = investment.get_service_charge()`nexpected = 0.50"
$ws.Range("F13").Value = "investment = InvestingAccount(12345, 123, 1000, date(2018, 5, 1), 3.00)`nexpected =`n            Account number: 12345`n            Balance: `$1,000.00`n            Date Created: 2018-05-01`n           Management Fee: `$3.00`n            Account Type: Investing"
$ws.Range("F12").Value = "investment = InvestingAccount(12345, 123, 1000, date(2013, 5, 1), 3.00)`nexpected =`n            Account number: 12345`n            Balance: `$1,000.00`n            Date Created: 2013-05-01`n           Management Fee: Waived`n            Account Type: Investing"
$ws.Range("G7").Value = "Account Number: 12345`nClient Number: 123`nBalance: 1000`nDate Created: 2025-5-1`nManagement Fee: 3.00"
$ws.Range("G8").Value = "Management Fee: 2.55"
$ws.Range("G9").Value = "Get Service Charge: 0.50"
$ws.Range("G10").Value = "Get Service Charge: 0.50"
$ws.Range("G11").Value = "Get Service Charge: 3.50"
$ws.Range("G13").Value = "Management fee: `$3.00`nMessage:`n            Account number: 12345`n            Balance: `$1,000.00`n            Date Created: 2018-05-01`n           Management Fee: `$3.00`n            Account Type: Investing"
$ws.Range("G12").Value = "Management Fee: Waived`nMessage:`n            Account number: 12345`n            Balance: `$1,000.00`n            Date Created: 2013-05-01`n           Management Fee: Waived`n            Account Type: Investing"

# Resize rows to fit the newly entered multi-line content
$ws.Rows.Item(7).RowHeight = 77.4
$ws.Rows.Item(8).RowHeight = 58.8
$ws.Rows.Item(9).RowHeight = 92.4
$ws.Rows.Item(10).RowHeight = 90
$ws.Rows.Item(11).RowHeight = 90
$ws.Rows.Item(12).RowHeight = 135
$ws.Rows.Item(13).RowHeight = 135.6

# Update the active selection to reflect where the author left off
$ws.Range("G12").Select() | Out-Null
